$wb = $excel.ActiveWorkbook

# "Generate Report for handoff": the 2513cad1 file (row 4 of each language
# sheet, status "Ready for handoff") has just had a new handoff generated,
# so its "Latest Handoff Datetime" (column D) is refreshed on both the
# zh-cn and de-de report sheets.

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-01-13 01:40:10"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-01-13 01:40:30"
